$d = $word.ActiveDocument

# "MINISTERIO DA FAZENDA" -> "MINISTERIO DA ECONOMIA"
$r1 = $d.Content
$r1.Find.Execute("FAZENDA", $false, $false, $false, $false, $false, $true, 1, $false, "ECONOMIA", 2)

# Consolidate the "OVR - " / "Operação" / " nº  {ovr_id}" runs into a single run
# (text is unchanged, but re-applying the replacement lets the engine merge the
# adjacent same-formatted runs, matching how Word re-serializes them on save).
$r2 = $d.Content
$r2.Find.Execute("Operação", $false, $false, $false, $false, $false, $true, 1, $false, "Operação", 2)
